# Updated data for Apr 7
#
# Column I holds "US Actual Confirmed Cases". Rows with real reported data use
# style s="14" (light fill) and a literal number; rows still being forecast use
# style s="15" (orange fill) with a shared growth-rate formula anchored on the
# first forecast row, i.e. "=<prevI>*(1+AVERAGE(<3 prior M cells>))".
#
# New actual data arrived for 4/7 (row 41), so:
#   - I41 (4/7) becomes a literal actual value -> style s="14"
#   - I40 (4/6) is no longer the first forecast day; it becomes a one-off,
#     manually overridden estimate (=I39*1.08) highlighted with the
#     yellow "note" style s="6"
#   - I42 (4/8) becomes the new first forecast day, anchoring the shared
#     growth formula based on the now-actual I41

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- I40: one-off manual estimate, highlighted (style s="6", like B3:B9) ---
$ws.Range("B3").Copy() | Out-Null
$ws.Range("I40").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("I40").Formula = "=I39*1.08"

# --- I41: new actual reported value for 4/7 (style s="14", like I39) ---
$ws.Range("I39").Copy() | Out-Null
$ws.Range("I41").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("I41").Value = 400335

# --- I42: becomes the new anchor of the forecast growth formula ---
$ws.Range("I42").Formula = "=I41*(1+AVERAGE(M39:M41))"

$excel.Calculate()

$ws.Range("A1").Select() | Out-Null
